# data: add 2023-04-14 notices
# Re-format the "公告明细" (notices detail) sheet:
#  - give the header row (A1:F1) the same "bold + thin border" look already
#    used by the "公告汇总" sheet's header, by copying its cell formatting
#  - strip the hyperlink styling/objects from the F column (links become
#    plain text)
#  - reset the page margins back to Excel's defaults

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 公告明细 (detail sheet, holds the new rows)
$ws2 = $wb.Worksheets.Item(2)   # 公告汇总 (summary sheet, already has the target header style)

# --- Header row formatting: match the style used on the summary sheet ---
$ws2.Range("A1:C1").Copy()
$ws1.Range("A1:F1").PasteSpecial(-4122)   # xlPasteFormats

# --- Drop the hyperlinks in column F (and their "Hyperlink" cell style) ---
$ws1.Range("F2:F8").Style = "Normal"
$ws1.Hyperlinks.Delete()

# --- Reset page margins to Excel defaults (inches -> points) ---
$ps = $ws1.PageSetup
$ps.LeftMargin   = 54    # 0.75"
$ps.RightMargin  = 54    # 0.75"
$ps.TopMargin    = 72    # 1"
$ps.BottomMargin = 72    # 1"
$ps.HeaderMargin = 36    # 0.5"
$ps.FooterMargin = 36    # 0.5"
